$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Add a new "Country" column (P) with "UK" for every recipe row (2-106)
$ws.Range("P1").Value = "Country"
$lastRow = 106
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 16).Value = "UK"
}

# 2. Fix typo in Saffron risotto ingredients (row 46, column B):
#    remove the stray word "tablespoon" from "15 ml tablespoon regular olive oil"
$fixed = "1 litre chicken stock`n0.4g saffron threads`n65 g butter`n15 ml regular olive oil`n50 g shallots`n250 g risotto rice`n125 ml dry marsala`n30 g parmesan`nsalt and pepper, to taste"
$ws.Range("B46").Value = $fixed

# Leave the selection on the cell that was edited, matching where the author
# was working when the workbook was last saved.
$ws.Range("C46").Select()
